$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (rows 386-464): date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newData = @(
    @(386, 44460, 1, 17, 94.56001779953276),
    @(387, 44461, 1, 18, 100.1223717877406),
    @(388, 44462, 2, 18, 100.1223717877406),
    @(389, 44463, 3, 13, 72.31060184670152),
    @(390, 44464, 1, 12, 66.74824785849371),
    @(391, 44465, 0, 8, 44.49883190566248),
    @(392, 44466, 2, 10, 55.6235398820781),
    @(393, 44467, 0, 9, 50.06118589387028),
    @(394, 44468, 0, 8, 44.49883190566248),
    @(395, 44469, 4, 10, 55.6235398820781),
    @(396, 44470, 2, 9, 50.06118589387028),
    @(397, 44471, 1, 9, 50.06118589387028),
    @(398, 44472, 2, 11, 61.1858938702859),
    @(399, 44473, 1, 10, 55.6235398820781),
    @(400, 44474, 0, 10, 55.6235398820781),
    @(401, 44475, 1, 11, 61.1858938702859),
    @(402, 44476, 2, 9, 50.06118589387028),
    @(403, 44477, 1, 8, 44.49883190566248),
    @(404, 44478, 0, 7, 38.93647791745467),
    @(405, 44479, 1, 6, 33.37412392924686),
    @(406, 44480, 0, 5, 27.81176994103905),
    @(407, 44481, 0, 5, 27.81176994103905),
    @(408, 44482, 0, 4, 22.24941595283124),
    @(409, 44483, 1, 3, 16.68706196462343),
    @(410, 44484, 4, 6, 33.37412392924686),
    @(411, 44485, 0, 6, 33.37412392924686),
    @(412, 44486, 1, 6, 33.37412392924686),
    @(413, 44487, 0, 6, 33.37412392924686),
    @(414, 44488, 1, 7, 38.93647791745467),
    @(415, 44489, 0, 7, 38.93647791745467),
    @(416, 44490, 3, 9, 50.06118589387028),
    @(417, 44491, 4, 9, 50.06118589387028),
    @(418, 44492, 3, 12, 66.74824785849371),
    @(419, 44493, 1, 12, 66.74824785849371),
    @(420, 44494, 1, 13, 72.31060184670152),
    @(421, 44495, 1, 13, 72.31060184670152),
    @(422, 44496, 1, 14, 77.87295583490933),
    @(423, 44497, 0, 11, 61.1858938702859),
    @(424, 44498, 0, 7, 38.93647791745467),
    @(425, 44499, 1, 5, 27.81176994103905),
    @(426, 44500, 0, 4, 22.24941595283124),
    @(427, 44501, 2, 5, 27.81176994103905),
    @(428, 44502, 1, 5, 27.81176994103905),
    @(429, 44503, 0, 4, 22.24941595283124),
    @(430, 44504, 0, 4, 22.24941595283124),
    @(431, 44505, 0, 4, 22.24941595283124),
    @(432, 44506, 0, 3, 16.68706196462343),
    @(433, 44507, 0, 3, 16.68706196462343),
    @(434, 44508, 0, 1, 5.56235398820781),
    @(435, 44509, 0, 0, 0),
    @(436, 44510, 0, 0, 0),
    @(437, 44511, 0, 0, 0),
    @(438, 44512, 0, 0, 0),
    @(439, 44513, 0, 0, 0),
    @(440, 44514, 2, 2, 11.12470797641562),
    @(441, 44515, 0, 2, 11.12470797641562),
    @(442, 44516, 17, 19, 105.6847257759484),
    @(443, 44517, 1, 20, 111.2470797641562),
    @(444, 44518, 0, 20, 111.2470797641562),
    @(445, 44519, 1, 21, 116.809433752364),
    @(446, 44520, 0, 21, 116.809433752364),
    @(447, 44521, 0, 19, 105.6847257759484),
    @(448, 44522, 0, 19, 105.6847257759484),
    @(449, 44523, 3, 5, 27.81176994103905),
    @(450, 44524, 28, 32, 177.9953276226499),
    @(451, 44525, 2, 34, 189.1200355990655),
    @(452, 44526, 1, 34, 189.1200355990655),
    @(453, 44527, 1, 35, 194.6823895872733),
    @(454, 44528, 1, 36, 200.2447435754811),
    @(455, 44529, 1, 37, 205.8070975636889),
    @(456, 44530, 0, 34, 189.1200355990655),
    @(457, 44531, 1, 7, 38.93647791745467),
    @(458, 44532, 1, 6, 33.37412392924686),
    @(459, 44533, 1, 6, 33.37412392924686),
    @(460, 44534, 2, 7, 38.93647791745467),
    @(461, 44535, 3, 9, 50.06118589387028),
    @(462, 44536, 1, 9, 50.06118589387028),
    @(463, 44537, 3, 12, 66.74824785849371),
    @(464, 44538, 0, 11, 61.1858938702859)
)

foreach ($entry in $newData) {
    $r = $entry[0]
    # Copy format (style) from the last existing data row (385) which carries the date style (s="2")
    $ws.Cells.Item(385, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}
